$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 23 ---
# Column C keeps its text ("Inv loss run...") but the shared-string slot it
# points at changes as a side effect of inserting new strings ahead of it;
# functionally the value is unchanged.
$ws.Range("C23").Value = "Inv loss run. MF-LRB enhanced repl target increased to 6%"
# Column E stays the same.
$ws.Range("E23").Value = "[10, 10, 7, 10]; [10, 10, 10, 10]"

# --- Add new row 24 (name/notes entered before the bracketed D values) ---
$ws.Range("A24").Value = "_complete_spectracomments_fixededp"
$ws.Range("B24").Value = 45961
$ws.Range("B24").NumberFormat = "d-mmm"
$ws.Range("C24").Value = 'Initial P-58 runs use "generate" mode: 1 EDP per run because it mirrors validation.'

# Column D's bracketed list gets the "-6.2" moved from the 3rd to the 2nd slot.
$ws.Range("D23").Value = "[-6, -6, -6, -6]; [-6, -6.2, -6, -7]"
$ws.Range("D24").Value = "[-6, -6, -6, -6]; [-6, -6.5, -6, -7]"
$ws.Range("E24").Value = "[10, 10, 7, 10]; [10, 10, 10, 10]"

# --- View state: scroll back to top-left and select the newly entered cell ---
$null = $ws.Range("A1").Select()
$null = $ws.Range("D24").Select()
